# Data clean-up: populate the new "MSI" column (D) in the lookup table
# located at rows 10-14 of Hoja1, and update the active selection to
# reflect where the user finished working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("D10").Value = "MSI"

# Data rows: fill in the MSI values that correspond to each SERIAL/MSIDN
$ws.Range("D11").Value = "732111198172294"
$ws.Range("D12").Value = "732111198172293"
$ws.Range("D13").Value = "732111198172294"
$ws.Range("D14").Value = "732111198172293"

# D11 previously held an (empty) hyperlink-style format; align it with the
# plain style already used by the rest of the data rows (same as C11).
$ws.Range("C11").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection where the user last clicked
$ws.Range("D19").Select()
